$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# --- ALC ---
$ws1.Range("H11").Value = 79.92856999999999
$ws1.Range("I11").Value = 79.92856999999999
$ws1.Range("K11").Value = 79.92856999999999
$ws1.Range("M11").Value = 60.07143000000001
$ws1.Range("H38").Value = 4772
$ws1.Range("I38").Value = 1053.4
$ws1.Range("J38").Value = 7428.143
$ws1.Range("K38").Value = 3160.2
$ws1.Range("L38").Value = 22284.429
$ws1.Range("M38").Value = -2788.2
$ws1.Range("N38").Value = -23028.429
$ws1.Range("H58").Value = 62500724
$ws1.Range("I58").Value = 62500724
$ws1.Range("J58").Value = 0
$ws1.Range("K58").Value = 187502172
$ws1.Range("L58").Value = 0
$ws1.Range("M58").Value = -187502022
$ws1.Range("N58").ClearContents()
$ws1.Range("H68").Value = 750018750
$ws1.Range("J68").Value = 750018750
$ws1.Range("L68").Value = 750018750
$ws1.Range("N68").Value = -750020248
$ws1.Range("H71").Value = 750018750
$ws1.Range("J71").Value = 750018750
$ws1.Range("L71").Value = 2250056250
$ws1.Range("N71").Value = -2250063738
$ws1.Range("H112").Value = 2732.3572
$ws1.Range("I112").Value = 1883
$ws1.Range("J112").Value = 2964
$ws1.Range("K112").Value = 5649
$ws1.Range("L112").Value = 8892
$ws1.Range("M112").Value = -4541
$ws1.Range("N112").Value = -11108
$ws1.Range("H132").Value = 1795.7858
$ws1.Range("I132").Value = 1126.5
$ws1.Range("K132").Value = 3379.5
$ws1.Range("M132").Value = -849.5
$ws1.Range("H138").Value = 4215.7915
$ws1.Range("I138").Value = 2373.5454
$ws1.Range("J138").Value = 5026.38
$ws1.Range("K138").Value = 7120.6362
$ws1.Range("L138").Value = 15079.14
$ws1.Range("M138").Value = -1980.6362
$ws1.Range("N138").Value = -25359.14

# --- ARM ---
$ws2.Range("H11").Value = 23333834
$ws2.Range("J11").Value = 0
$ws2.Range("L11").Value = 0
$ws2.Range("N11").ClearContents()
$ws2.Range("H32").Value = 5045.4595
$ws2.Range("I32").Value = 4664.5806
$ws2.Range("K32").Value = 4664.5806
$ws2.Range("M32").Value = -4377.5806
$ws2.Range("H45").Value = 3947.5
$ws2.Range("I45").Value = 1200
$ws2.Range("K45").Value = 1200
$ws2.Range("M45").Value = -823
$ws2.Range("H74").Value = 2983.111
$ws2.Range("I74").Value = 3106
$ws2.Range("K74").Value = 3106
$ws2.Range("M74").Value = -2232
$ws2.Range("H77").Value = 2983.111
$ws2.Range("I77").Value = 3106
$ws2.Range("K77").Value = 15530
$ws2.Range("M77").Value = -11162
$ws2.Range("H132").Value = 2220.8462
$ws2.Range("I132").Value = 2249.75
$ws2.Range("J132").Value = 1874
$ws2.Range("K132").Value = 6749.25
$ws2.Range("L132").Value = 5622
$ws2.Range("M132").Value = -4219.25
$ws2.Range("N132").Value = -10682

# --- BSM ---
$ws3.Range("H94").Value = 1159.125
$ws3.Range("I94").Value = 1167.5714
$ws3.Range("K94").Value = 1167.5714
$ws3.Range("M94").Value = -716.5714
$ws3.Range("H99").Value = 3528.2778
$ws3.Range("I99").Value = 1600.9
$ws3.Range("J99").Value = 5937.5
$ws3.Range("K99").Value = 1600.9
$ws3.Range("L99").Value = 5937.5
$ws3.Range("M99").Value = -102.9000000000001
$ws3.Range("N99").Value = -8933.5
$ws3.Range("H132").Value = 256570.28
$ws3.Range("J132").Value = 256570.28
$ws3.Range("L132").Value = 256570.28
$ws3.Range("N132").Value = -266690.28
$ws3.Range("H133").Value = 79997
$ws3.Range("J133").Value = 79997
$ws3.Range("L133").Value = 79997
$ws3.Range("N133").Value = -90117
$ws3.Range("H134").Value = 2292.4285
$ws3.Range("I134").Value = 2230.4644
$ws3.Range("J134").Value = 2540.2856
$ws3.Range("K134").Value = 6691.3932
$ws3.Range("L134").Value = 7620.8568
$ws3.Range("M134").Value = -4156.3932
$ws3.Range("N134").Value = -12690.8568

# --- CRP ---
$ws4.Range("H31").Value = 3633.426
$ws4.Range("I31").Value = 1541.35
$ws4.Range("J31").Value = 4864.0586
$ws4.Range("K31").Value = 1541.35
$ws4.Range("L31").Value = 4864.0586
$ws4.Range("M31").Value = -1246.35
$ws4.Range("N31").Value = -5454.0586
$ws4.Range("H34").Value = 3633.426
$ws4.Range("I34").Value = 1541.35
$ws4.Range("J34").Value = 4864.0586
$ws4.Range("K34").Value = 1541.35
$ws4.Range("L34").Value = 4864.0586
$ws4.Range("M34").Value = -1339.35
$ws4.Range("N34").Value = -5268.0586
$ws4.Range("H58").Value = 1620.5952
$ws4.Range("I58").Value = 1820.96
$ws4.Range("K58").Value = 1820.96
$ws4.Range("M58").Value = -1617.96
$ws4.Range("H94").Value = 2685.88
$ws4.Range("J94").Value = 3396.077
$ws4.Range("L94").Value = 3396.077
$ws4.Range("N94").Value = -4298.077
$ws4.Range("H132").Value = 1598.9722
$ws4.Range("I132").Value = 1628.3529
$ws4.Range("J132").Value = 1099.5
$ws4.Range("K132").Value = 4885.0587
$ws4.Range("L132").Value = 3298.5
$ws4.Range("M132").Value = -2355.0587
$ws4.Range("N132").Value = -8358.5
$ws4.Range("H134").Value = 2025.4857
$ws4.Range("I134").Value = 2003
$ws4.Range("K134").Value = 6009
$ws4.Range("M134").Value = -3474
$ws4.Range("H136").Value = 1620.5952
$ws4.Range("I136").Value = 1820.96
$ws4.Range("K136").Value = 5462.88
$ws4.Range("M136").Value = -2912.88

# --- CUL ---
$ws5.Range("H2").Value = 78.545456
$ws5.Range("I2").Value = 70.333336
$ws5.Range("J2").Value = 81.625
$ws5.Range("K2").Value = 422.000016
$ws5.Range("L2").Value = 489.75
$ws5.Range("M2").Value = -309.000016
$ws5.Range("N2").Value = -715.75
$ws5.Range("H23").Value = 5460
$ws5.Range("I23").Value = 20001
$ws5.Range("J23").Value = 3036.5
$ws5.Range("K23").Value = 60003
$ws5.Range("L23").Value = 9109.5
$ws5.Range("M23").Value = -59768
$ws5.Range("N23").Value = -9579.5
$ws5.Range("H41").Value = 622.875
$ws5.Range("I41").Value = 200
$ws5.Range("J41").Value = 763.8333
$ws5.Range("K41").Value = 600
$ws5.Range("L41").Value = 2291.4999
$ws5.Range("M41").Value = -262
$ws5.Range("N41").Value = -2967.4999
$ws5.Range("H63").Value = 2450
$ws5.Range("I63").Value = 2450
$ws5.Range("K63").Value = 7350
$ws5.Range("M63").Value = -6601
$ws5.Range("H66").Value = 2450
$ws5.Range("I66").Value = 2450
$ws5.Range("K66").Value = 22050
$ws5.Range("M66").Value = -18306
$ws5.Range("H107").Value = 1361.6842
$ws5.Range("I107").Value = 2201.5
$ws5.Range("K107").Value = 6604.5
$ws5.Range("M107").Value = -4684.5
$ws5.Range("H128").Value = 255833.17
$ws5.Range("I128").Value = 255833.17
$ws5.Range("K128").Value = 767499.51
$ws5.Range("M128").Value = -762519.51

# --- GSM ---
$ws6.Range("H59").Value = 53332.668
$ws6.Range("I59").Value = 52499.5
$ws6.Range("J59").Value = 54999
$ws6.Range("K59").Value = 52499.5
$ws6.Range("L59").Value = 54999
$ws6.Range("N59").Value = -56165
$ws6.Range("M59").Value = -51916.5
$ws6.Range("H97").Value = 341.94116
$ws6.Range("I97").Value = 296.8
$ws6.Range("K97").Value = 296.8
$ws6.Range("M97").Value = 199.2
$ws6.Range("H122").Value = 5910.933
$ws6.Range("J122").Value = 3304.2666
$ws6.Range("L122").Value = 9912.799800000001
$ws6.Range("N122").Value = -14812.7998
$ws6.Range("H132").Value = 2082.9048
$ws6.Range("I132").Value = 1915.7059
$ws6.Range("K132").Value = 5747.1177
$ws6.Range("M132").Value = -3217.1177

# --- LTW ---
$ws7.Range("H46").Value = 2803.25
$ws7.Range("I46").Value = 1099.25
$ws7.Range("J46").Value = 3371.25
$ws7.Range("K46").Value = 1099.25
$ws7.Range("L46").Value = 3371.25
$ws7.Range("M46").Value = -911.25
$ws7.Range("N46").Value = -3747.25
$ws7.Range("H55").Value = 1861.8422
$ws7.Range("I55").Value = 307.81818
$ws7.Range("K55").Value = 307.81818
$ws7.Range("M55").Value = -134.81818
$ws7.Range("H69").Value = 750012500
$ws7.Range("J69").Value = 750012500
$ws7.Range("L69").Value = 750012500
$ws7.Range("N69").Value = -750014122
$ws7.Range("H72").Value = 750012500
$ws7.Range("J72").Value = 750012500
$ws7.Range("L72").Value = 2250037500
$ws7.Range("N72").Value = -2250045612
$ws7.Range("H93").Value = 2311.2083
$ws7.Range("I93").Value = 791.6667
$ws7.Range("K93").Value = 791.6667
$ws7.Range("M93").Value = 456.3333
$ws7.Range("H100").Value = 7045.364
$ws7.Range("I100").Value = 4500
$ws7.Range("J100").Value = 7999.875
$ws7.Range("K100").Value = 4500
$ws7.Range("L100").Value = 7999.875
$ws7.Range("M100").Value = -3959
$ws7.Range("N100").Value = -9081.875
$ws7.Range("H122").Value = 8565.593999999999
$ws7.Range("I122").Value = 6519.636
$ws7.Range("J122").Value = 13066.7
$ws7.Range("K122").Value = 19558.908
$ws7.Range("L122").Value = 39200.10000000001
$ws7.Range("M122").Value = -17108.908
$ws7.Range("N122").Value = -44100.10000000001
$ws7.Range("H132").Value = 2118.8728
$ws7.Range("I132").Value = 1939.875
$ws7.Range("K132").Value = 5819.625
$ws7.Range("M132").Value = -3289.625

# --- WVR ---
$ws8.Range("H113").Value = 726.88
$ws8.Range("I113").Value = 409.0909
$ws8.Range("K113").Value = 1227.2727
$ws8.Range("M113").Value = 942.7273
$ws8.Range("H126").Value = 2237.9443
$ws8.Range("I126").Value = 1900.4231
$ws8.Range("J126").Value = 3115.5
$ws8.Range("K126").Value = 5701.2693
$ws8.Range("L126").Value = 9346.5
$ws8.Range("M126").Value = -3231.2693
$ws8.Range("N126").Value = -14286.5
$ws8.Range("H136").Value = 1943.174
$ws8.Range("J136").Value = 2875.125
$ws8.Range("L136").Value = 8625.375
$ws8.Range("N136").Value = -13725.375
